$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Update the text in E20 (shared string "XML + Add Advancement" -> appended text)
$ws.Range("E20").Value = "XML + Add Advancement, Feats being added with partial filter"

# Make the E20 cell wrap text (matches style used by other long comment rows)
$ws.Range("E20").WrapText = $true

# Update end time in C20, which cascades through the dependent formulas in D/F/G columns
$ws.Range("C20").Value = 0.98958333333333337

# Set the row height for row 20 to accommodate wrapped text
$ws.Rows.Item(20).RowHeight = 30

# Update the active selection to E13 as recorded in the saved view state
$ws.Range("E13").Select()
